# Apply updated "dSF" (column F) values for rows 2-12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = 14
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -1
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -8
$ws.Range("F12").Value = 2
